$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Tuesday", "ragu", 0.6830000000000001, 6, 0.6830000000000001, '[{"name": "Diet Base", "total": {"protein": 0, "calories": 0, "cost": 0.0, "amount": 0.0, "name": []}}]'),
    @("Tuesday", "ragu", 0.6830000000000001, 6, 0.6830000000000001, '[{"name": "Diet Base", "total": {"protein": 0, "calories": 0, "cost": 0.0, "amount": 0.0, "name": []}}]'),
    @("Tuesday", "ragu", 3.08716, 27.12, 3.08716, '[{"name": "Diet Base", "total": {"protein": 0, "calories": 0, "cost": 0.0, "amount": 0.0, "name": []}}]'),
    @("Tuesday", "ragu", 1.92606, 16.92, 1.92606, '[{"name": "Diet Base", "total": {"protein": 0, "calories": 0, "cost": 0.0, "amount": 0.0, "name": []}}]')
)

$startRow = 5
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
